$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $orig = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $orig
}

# Row 2
Set-TextValue 2 4 '29.173.63'
Set-TextValue 2 5 '  -0.80%  '

# Row 3
Set-TextValue 3 4 '1.858.59'
Set-TextValue 3 5 '  -1.07%  '

# Row 4
Set-TextValue 4 4 '0.9998'
Set-TextValue 4 5 '  -0.12%  '

# Row 5
Set-TextValue 5 4 '242.14'
Set-TextValue 5 5 '  -0.66%  '

# Row 6
Set-TextValue 6 4 '0.7012'
Set-TextValue 6 5 '  -2.27%  '

# Row 7
Set-TextValue 7 4 '1.000'
Set-TextValue 7 5 '  -0.09%  '

# Row 8
Set-TextValue 8 4 '0.07797'
Set-TextValue 8 5 '  -2.22%  '

# Row 9
Set-TextValue 9 4 '0.3105'
Set-TextValue 9 5 '  -1.40%  '

# Row 10
Set-TextValue 10 4 '23.88'
Set-TextValue 10 5 '  -4.27%  '

# Row 11
Set-TextValue 11 4 '0.07798'
Set-TextValue 11 5 '  -3.99%  '

# Row 12
Set-TextValue 12 4 '1.862.07'
Set-TextValue 12 5 '  -1.28%  '

# Row 13
Set-TextValue 13 2 'Polkadot'
Set-TextValue 13 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 13 4 '5.125'
Set-TextValue 13 5 '  -2.09%  '

# Row 14
Set-TextValue 14 2 'Litecoin'
Set-TextValue 14 3 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 14 4 '92.60'
Set-TextValue 14 5 '  -2.19%  '

# Row 15
Set-TextValue 15 4 '0.6904'
Set-TextValue 15 5 '  -2.45%  '

# Row 16
Set-TextValue 16 4 '6.554'
Set-TextValue 16 5 '  +2.26%  '

# Row 17
Set-TextValue 17 4 '0.000008447'
Set-TextValue 17 5 '  +0.17%  '

# Row 18
Set-TextValue 18 4 '29.213.24'
Set-TextValue 18 5 '  -0.68%  '

# Row 19
Set-TextValue 19 4 '250.01'
Set-TextValue 19 5 '  -1.26%  '

# Row 20
Set-TextValue 20 4 '2.111.52'
Set-TextValue 20 5 '  -1.09%  '

# Row 21
Set-TextValue 21 4 '12.92'
Set-TextValue 21 5 '  -3.24%  '

# Row 22
Set-TextValue 22 4 '0.9998'
Set-TextValue 22 5 '  -0.07%  '

# Row 23
Set-TextValue 23 4 '7.603'
Set-TextValue 23 5 '  -0.97%  '

# Row 24
Set-TextValue 24 4 '1.000'
Set-TextValue 24 5 '  -0.15%  '

# Row 25
Set-TextValue 25 4 '0.1532'
Set-TextValue 25 5 '  -3.23%  '

# Row 26
Set-TextValue 26 4 '160.72'
Set-TextValue 26 5 '  -0.85%  '

# Row 27
Set-TextValue 27 4 '8.908'
Set-TextValue 27 5 '  -1.82%  '

# Row 28
Set-TextValue 28 4 '18.57'
Set-TextValue 28 5 '  -2.09%  '

# Row 29
Set-TextValue 29 4 '1.570'
Set-TextValue 29 5 '  +4.08%  '

# Row 30
Set-TextValue 30 4 '4.272'
Set-TextValue 30 5 '  -3.32%  '

# Row 31
Set-TextValue 31 4 '4.252'
Set-TextValue 31 5 '  -1.50%  '

# Row 32
Set-TextValue 32 4 '1.211'
Set-TextValue 32 5 '  -1.12%  '

# Row 33
Set-TextValue 33 4 '0.05235'
Set-TextValue 33 5 '  -1.55%  '

# Row 34
Set-TextValue 34 4 '0.7576'
Set-TextValue 34 5 '  -0.15%  '

# Row 35
Set-TextValue 35 4 '1.876'
Set-TextValue 35 5 '  -3.75%  '

# Row 36
Set-TextValue 36 4 '1.176'
Set-TextValue 36 5 '  +0.04%  '

# Row 37
Set-TextValue 37 4 '2.708'
Set-TextValue 37 5 '  +0.09%  '

# Row 38
Set-TextValue 38 4 '0.01862'
Set-TextValue 38 5 '  -1.58%  '

# Row 39
Set-TextValue 39 4 '1.221.76'
Set-TextValue 39 5 '  -4.01%  '

# Row 40
Set-TextValue 40 4 '2.721'
Set-TextValue 40 5 '  -1.54%  '

# Row 41
Set-TextValue 41 4 '0.8996'
Set-TextValue 41 5 '  -0.69%  '

# Row 42
Set-TextValue 42 4 '110.49'
Set-TextValue 42 5 '  -1.23%  '

# Row 43
Set-TextValue 43 4 '5.826'
Set-TextValue 43 5 '  -9.04%  '

# Row 44
Set-TextValue 44 4 '0.9997'
Set-TextValue 44 5 '  -0.11%  '

# Row 45
Set-TextValue 45 4 '67.39'
Set-TextValue 45 5 '  -9.38%  '

# Row 46
Set-TextValue 46 4 '2.008.70'
Set-TextValue 46 5 '  -1.03%  '

# Row 47
Set-TextValue 47 2 'Mantle'
Set-TextValue 47 3 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 47 4 '0.5182'
Set-TextValue 47 5 '  -0.48%  '

# Row 48
Set-TextValue 48 2 'EnergySwap'
Set-TextValue 48 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 48 4 '9.538'
Set-TextValue 48 5 '  +0.04%  '

# Row 49
Set-TextValue 49 2 'BabyDogeCoin'
Set-TextValue 49 3 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 49 4 '0.00000000121'
Set-TextValue 49 5 '  -5.96%  '

# Row 50
Set-TextValue 50 4 '1.767'
Set-TextValue 50 5 '  -2.25%  '

# Row 51
Set-TextValue 51 4 '7.031'
Set-TextValue 51 5 '  -1.02%  '
